$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 swap
$ws.Range("N2").Value = 1
$ws.Range("Q2").Value = 0

# Row 4 updates
$ws.Range("N4").Value = 0.523252976771423
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0.04650595354284603
$ws.Range("Q4").Value = 0.7361990430222111
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = 0.4723980860444221
$ws.Range("T4").Value = 0.7699568058175457
$ws.Range("V4").Value = 0.5399136116350913

# Column B updates (rows 15-21)
$ws.Range("B15").Value = 10.14140507033079
$ws.Range("B16").Value = 4.865859824874436
$ws.Range("B17").Value = -5.030967384128449
$ws.Range("B18").Value = 6.862334569795242
$ws.Range("B19").Value = 10.03672408650753
$ws.Range("B20").Value = 13.40550741141245
$ws.Range("B21").Value = 24.810071406628
